# Add a new "Chromium(VI)" CFR method row to the CFR Methods worksheet.
# This inserts a new row at position 75 (pushing existing rows 75-180 down
# to 76-181) and populates it with the Chromium(VI) method data, matching
# the pattern of the surrounding rows (Char_Name, CASNumber, Method_Code,
# Method_Context, CFR_Method).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 75; Excel shifts rows 75:180 down to 76:181
# and copies formatting from the row above, same as a normal Excel insert.
$ws.Rows("75:75").Insert()

# Populate the new row. Write A (Char_Name) first, then C (Method_Code)
# before B (CASNumber) so the new shared-string entries land in the same
# order as the source workbook (Chromium(VI), 218.6, 18540-29-9).
$ws.Range("A75").Value = "Chromium(VI)"
$ws.Range("C75").Value = "218.6"
$ws.Range("B75").Value = "18540-29-9"
$ws.Range("D75").Value = "U.S. Environmental Protection Agency"
$ws.Range("E75").Value = "Yes"

# The B column (CASNumber) on the new row carries no explicit style, unlike
# the inherited style copied from the row above during the insert.
$ws.Range("B75").Style = "Normal"

# Match the author's final selection/view state on the sheet.
$ws.Range("B75").Select()
